$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# S05/G01 status rows (34-36): mark tasks implemented and fill in
# deviations / remarks / pending-work notes for the Zerodha Kite Connect
# client integration.

$ws.Range("G34").Value = "implemented"
$ws.Range("F34").Value = "Introduced ZerodhaClient wrapper reading API key from Settings and lazily instantiating KiteConnect."
$ws.Range("H34").Value = "Backend can now construct a Zerodha client once an access token is available."
$ws.Range("I34").Value = "Hook this client into OAuth/token storage and execution flows in S05/G02–G03."

$ws.Range("G35").Value = "implemented"
$ws.Range("F35").Value = "Implemented ZerodhaClient.place_order that composes KiteConnect place_order parameters."
$ws.Range("H35").Value = "Order placement service is tested against a fake Kite client (no real network calls)."
$ws.Range("I35").Value = "Wire real order placement into manual queue execution once broker connection is ready."

$ws.Range("G36").Value = "implemented"
$ws.Range("F36").Value = "Exposed ZerodhaClient.list_orders and get_order_history as thin wrappers over KiteConnect APIs."
$ws.Range("H36").Value = "Backend has basic services to retrieve Zerodha order book and order history."
$ws.Range("I36").Value = "Use these services for status sync and richer order views in S07."
